$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 499.6
$ws.Range("I2").Value = 488.25
$ws.Range("K2").Value = 488.25
$ws.Range("M2").Value = -375.25

$ws.Range("H75").Value = 209795.36
$ws.Range("J75").Value = 209795.36
$ws.Range("L75").Value = 209795.36
$ws.Range("N75").Value = -211667.36

$ws.Range("H78").Value = 209795.36
$ws.Range("J78").Value = 209795.36
$ws.Range("L78").Value = 629386.08
$ws.Range("N78").Value = -638746.08

$ws.Range("H107").Value = 6951.643
$ws.Range("I107").Value = 8393.091
$ws.Range("J107").Value = 1666.3334
$ws.Range("K107").Value = 8393.091
$ws.Range("L107").Value = 1666.3334
$ws.Range("M107").Value = -6473.091
$ws.Range("N107").Value = -5506.3334

$ws.Range("H112").Value = 2754.4285
$ws.Range("J112").Value = 2747.8125
$ws.Range("L112").Value = 8243.4375
$ws.Range("N112").Value = -10459.4375

$ws.Range("H128").Value = 89993
$ws.Range("J128").Value = 89993
$ws.Range("L128").Value = 89993
$ws.Range("N128").Value = -99953

$ws.Range("H132").Value = 1697491.9
$ws.Range("I132").Value = 2444.0728
$ws.Range("K132").Value = 7332.2184
$ws.Range("M132").Value = -4802.2184

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16479.535
$ws.Range("I32").Value = 14843.077
$ws.Range("K32").Value = 14843.077
$ws.Range("M32").Value = -14556.077

$ws.Range("H102").Value = 10366.086
$ws.Range("I102").Value = 12014.904
$ws.Range("J102").Value = 7892.857
$ws.Range("K102").Value = 12014.904
$ws.Range("L102").Value = 7892.857
$ws.Range("M102").Value = -10392.904
$ws.Range("N102").Value = -11136.857

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5013
$ws.Range("J20").Value = 7913.6665
$ws.Range("L20").Value = 7913.6665
$ws.Range("N20").Value = -8407.6665

$ws.Range("H105").Value = 43836.926
$ws.Range("I105").Value = 159371.28
$ws.Range("K105").Value = 159371.28
$ws.Range("M105").Value = -157624.28

$ws.Range("H130").Value = 87499.25
$ws.Range("J130").Value = 87499.25
$ws.Range("L130").Value = 87499.25
$ws.Range("N130").Value = -97539.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 657.0714
$ws.Range("J22").Value = 943.75
$ws.Range("L22").Value = 943.75
$ws.Range("N22").Value = -1643.75

$ws.Range("H99").Value = 8339095
$ws.Range("I99").Value = 20837330
$ws.Range("K99").Value = 20837330
$ws.Range("M99").Value = -20835832

$ws.Range("H122").Value = 10255.866
$ws.Range("I122").Value = 11126.077
$ws.Range("K122").Value = 33378.231
$ws.Range("M122").Value = -30928.231

$ws.Range("H126").Value = 8339095
$ws.Range("I126").Value = 20837330
$ws.Range("K126").Value = 62511990
$ws.Range("M126").Value = -62509520

$ws.Range("H141").Value = 595125.7
$ws.Range("J141").Value = 630491.8
$ws.Range("L141").Value = 630491.8
$ws.Range("N141").Value = -640851.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 287705.06
$ws.Range("I5").Value = 944.8125
$ws.Range("J5").Value = 529187.4
$ws.Range("K5").Value = 2834.4375
$ws.Range("L5").Value = 1587562.2
$ws.Range("M5").Value = -2722.4375
$ws.Range("N5").Value = -1587786.2

$ws.Range("H87").Value = 12865.467
$ws.Range("I87").Value = 7257.3125
$ws.Range("J87").Value = 19274.785
$ws.Range("K87").Value = 21771.9375
$ws.Range("L87").Value = 57824.355
$ws.Range("M87").Value = -20523.9375
$ws.Range("N87").Value = -60320.355

$ws.Range("H90").Value = 12865.467
$ws.Range("I90").Value = 7257.3125
$ws.Range("J90").Value = 19274.785
$ws.Range("K90").Value = 65315.8125
$ws.Range("L90").Value = 173473.065
$ws.Range("M90").Value = -59075.8125
$ws.Range("N90").Value = -185953.065

$ws.Range("H122").Value = 5458.6216
$ws.Range("I122").Value = 1306
$ws.Range("J122").Value = 6996.6294
$ws.Range("K122").Value = 11754
$ws.Range("L122").Value = 62969.6646
$ws.Range("M122").Value = -9304
$ws.Range("N122").Value = -67869.66459999999

$ws.Range("H135").Value = 287705.06
$ws.Range("I135").Value = 944.8125
$ws.Range("J135").Value = 529187.4
$ws.Range("K135").Value = 8503.3125
$ws.Range("L135").Value = 4762686.600000001
$ws.Range("M135").Value = -5968.3125
$ws.Range("N135").Value = -4767756.600000001

$ws.Range("H138").Value = 1608.8889
$ws.Range("I138").Value = 1622.5
$ws.Range("J138").Value = 1500
$ws.Range("K138").Value = 4867.5
$ws.Range("L138").Value = 4500
$ws.Range("M138").Value = 272.5
$ws.Range("N138").Value = -14780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1384.6666
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H102").Value = 4188.2554
$ws.Range("I102").Value = 4497.025
$ws.Range("J102").Value = 2423.8572
$ws.Range("K102").Value = 4497.025
$ws.Range("L102").Value = 2423.8572
$ws.Range("M102").Value = -2875.025
$ws.Range("N102").Value = -5667.8572

$ws.Range("H123").Value = 29794.834
$ws.Range("J123").Value = 29794.834
$ws.Range("L123").Value = 29794.834
$ws.Range("N123").Value = -34694.834

$ws.Range("H132").Value = 3814.2444
$ws.Range("I132").Value = 3920.024
$ws.Range("K132").Value = 11760.072
$ws.Range("M132").Value = -9230.072

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4615.5
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 5923.25
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 5923.25
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -7421.25

$ws.Range("H71").Value = 4615.5
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 5923.25
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 29616.25
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -37104.25

$ws.Range("H93").Value = 5372.5
$ws.Range("J93").Value = 1599.3334
$ws.Range("L93").Value = 1599.3334
$ws.Range("N93").Value = -4095.3334

$ws.Range("H100").Value = 5833.4443
$ws.Range("J100").Value = 8000
$ws.Range("L100").Value = 8000
$ws.Range("N100").Value = -9082
